$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-16 Tuesday", "2025-09-17 Wednesday"),
    @("89×67=", "33×86="),
    @("16×68=", "93×20="),
    @("80×27=", "35×58="),
    @("24×25=", "66×20="),
    @("87×62=", "60×88="),
    @("29×19=", "80×45="),
    @("47×67=", "87×27="),
    @("20×85=", "29×68="),
    @("96×63=", "46×98="),
    @("11×24=", "49×57="),
    @("62×63=", "82×94="),
    @("66×93=", "44×18="),
    @("45×31=", "15×28="),
    @("77×39=", "11×49="),
    @("39×36=", "48×20="),
    @("84×77=", "28×82="),
    @("50×93=", "79×70="),
    @("83×11=", "32×45="),
    @("27×34=", "68×99="),
    @("20×40=", "82×13="),
    @("23×16=", "19×16="),
    @("15×50=", "56×75="),
    @("34×34=", "51×87="),
    @("78×98=", "36×28="),
    @("34×98=", "60×61=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
